# Build the Dutch "Beschikbaarheid CO2 vervloeiing" availability template
# into Sheet1 of the workbook, matching the target layout/styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column widths ----
# Target XML widths are 37.44140625 (col A) and 10.33203125 (col C).
# This engine snaps ColumnWidth (chars) onto a 1/6-character grid when
# writing the XML <col> width, so we pick inputs whose rounded result is
# the nearest achievable grid point to the authored values.
$ws.Columns.Item(1).ColumnWidth = 36.65      # -> xml width 37.5  (closest to 37.44140625)
$ws.Columns.Item(3).ColumnWidth = 9.5        # -> xml width 10.333333333333334 (closest to 10.33203125)

# ---- Helper color for the light-gray "##" fill cells (theme Background1, -15% tint) ----
$grayFill = 14277081   # RGB 0xD9D9D9

# ---- Row 1 : title row ----
$a1 = $ws.Range("A1")
$a1.Value = "Beschikbaarheid CO2 vervloeiing"
$a1.Font.Bold = $true
$a1.Font.Underline = $true
$a1.Borders.Item(7).LineStyle = 1   # left
$a1.Borders.Item(8).LineStyle = 1   # top

$b1 = $ws.Range("B1")
$b1.Borders.Item(8).LineStyle = 1   # top

$c1 = $ws.Range("C1")
$c1.Borders.Item(8).LineStyle = 1   # top
$c1.Borders.Item(10).LineStyle = 1  # right

# ---- Rows 2-8 : label / ## / unit rows ----
$labels = @(
  "aantal draaiuren",
  "standby",
  "storing (verantwoordelijkheid Bright)",
  "storing (verantwoordelijkheid klant)",
  "stilstand t.g.v. gepland onderhoud",
  "stilstand t.g.v. niet-gepland onderhoud",
  "beschikbaarheid maand"
)
$units = @("uur", "uur", "uur", "uur", "uur", "uur", "%")

for ($i = 0; $i -lt $labels.Length; $i++) {
  $r = $i + 2
  $aCell = $ws.Cells.Item($r, 1)
  $aCell.Value = $labels[$i]
  $aCell.Borders.Item(7).LineStyle = 1   # left

  $bCell = $ws.Cells.Item($r, 2)
  $bCell.Value = "##"
  $bCell.Interior.Color = $grayFill

  $cCell = $ws.Cells.Item($r, 3)
  $cCell.Value = $units[$i]
  $cCell.Borders.Item(10).LineStyle = 1  # right
}

# ---- Row 9 : blank spacer row (keeps left/right borders only, no fill) ----
$a9 = $ws.Range("A9")
$a9.Borders.Item(7).LineStyle = 1   # left

$c9 = $ws.Range("C9")
$c9.Borders.Item(10).LineStyle = 1  # right

# ---- Row 10 : rolling-year availability row ----
$a10 = $ws.Range("A10")
$a10.Value = "beschikbaarheid voortschrijdend kalenderjaar (garantie 97 %)"
$a10.WrapText = $true
$a10.Borders.Item(7).LineStyle = 1   # left
$a10.Borders.Item(9).LineStyle = 1   # bottom

$b10 = $ws.Range("B10")
$b10.Value = "##"
$b10.Interior.Color = $grayFill
$b10.Borders.Item(9).LineStyle = 1   # bottom

$c10 = $ws.Range("C10")
$c10.Value = "%"
$c10.Borders.Item(10).LineStyle = 1  # right
$c10.Borders.Item(9).LineStyle = 1   # bottom

$ws.Rows.Item(10).RowHeight = 28.8

# ---- Selection / view ----
$ws.Range("A1:C10").Select()

Write-Host "template written"
